$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 20275.732
$ws.Range("I32").Value = 2406.4314
$ws.Range("J32").Value = 202542.6
$ws.Range("K32").Value = 2406.4314
$ws.Range("L32").Value = 202542.6
$ws.Range("M32").Value = -2119.4314
$ws.Range("N32").Value = -203116.6

# Row 92 (ARM)
$ws.Range("H92").Value = 17266.334
$ws.Range("J92").Value = 17266.334
$ws.Range("L92").Value = 17266.334
$ws.Range("N92").Value = -22258.334

# Row 97 (ARM)
$ws.Range("H97").Value = 13893424
$ws.Range("I97").Value = 15156427
$ws.Range("J97").Value = 390
$ws.Range("K97").Value = 15156427
$ws.Range("L97").Value = 390
$ws.Range("M97").Value = -15155931
$ws.Range("N97").Value = -1382

# Row 112 (ARM)
$ws.Range("H112").Value = 17290.25
$ws.Range("J112").Value = 17290.25
$ws.Range("L112").Value = 17290.25
$ws.Range("N112").Value = -20244.25

# Row 130 (ARM)
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Row 137 (ARM)
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

$ws = $wb.Worksheets.Item("BSM")
# Row 57 (BSM)
$ws.Range("H57").Value = 100500
$ws.Range("J57").Value = 100500
$ws.Range("L57").Value = 100500
$ws.Range("N57").Value = -101940

# Row 58 (BSM)
$ws.Range("H58").Value = 37250
$ws.Range("J58").Value = 37250
$ws.Range("L58").Value = 37250
$ws.Range("N58").Value = -37838

# Row 99 (BSM)
$ws.Range("H99").Value = 2354.7334
$ws.Range("I99").Value = 2101.111
$ws.Range("J99").Value = 2735.1667
$ws.Range("K99").Value = 2101.111
$ws.Range("L99").Value = 2735.1667
$ws.Range("M99").Value = -603.1109999999999
$ws.Range("N99").Value = -5731.1667

# Row 103 (BSM)
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# Row 131 (BSM)
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# Row 136 (BSM)
$ws.Range("H136").Value = 100500
$ws.Range("J136").Value = 100500
$ws.Range("L136").Value = 100500
$ws.Range("N136").Value = -110700

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 1931.8462
$ws.Range("I31").Value = 1131.875
$ws.Range("J31").Value = 3211.8
$ws.Range("K31").Value = 1131.875
$ws.Range("L31").Value = 3211.8
$ws.Range("M31").Value = -836.875
$ws.Range("N31").Value = -3801.8

# Row 34 (CRP)
$ws.Range("H34").Value = 1931.8462
$ws.Range("I34").Value = 1131.875
$ws.Range("J34").Value = 3211.8
$ws.Range("K34").Value = 1131.875
$ws.Range("L34").Value = 3211.8
$ws.Range("M34").Value = -929.875
$ws.Range("N34").Value = -3615.8

# Row 47 (CRP)
$ws.Range("H47").Value = 4925.7144
$ws.Range("I47").Value = 3990
$ws.Range("K47").Value = 3990
$ws.Range("M47").Value = -3424

# Row 58 (CRP)
$ws.Range("H58").Value = 2001.4524
$ws.Range("I58").Value = 768.6799999999999
$ws.Range("J58").Value = 3814.353
$ws.Range("K58").Value = 768.6799999999999
$ws.Range("L58").Value = 3814.353
$ws.Range("M58").Value = -565.6799999999999
$ws.Range("N58").Value = -4220.353

# Row 107 (CRP)
$ws.Range("H107").Value = 330.75
$ws.Range("J107").Value = 348.64285
$ws.Range("L107").Value = 348.64285
$ws.Range("N107").Value = -4188.64285

# Row 115 (CRP)
$ws.Range("H115").Value = 24853.334
$ws.Range("J115").Value = 24853.334
$ws.Range("L115").Value = 24853.334
$ws.Range("N115").Value = -27203.334

# Row 122 (CRP)
$ws.Range("H122").Value = 1686.2593
$ws.Range("I122").Value = 1133.4375
$ws.Range("J122").Value = 2490.3635
$ws.Range("K122").Value = 3400.3125
$ws.Range("L122").Value = 7471.0905
$ws.Range("M122").Value = -950.3125
$ws.Range("N122").Value = -12371.0905

# Row 132 (CRP)
$ws.Range("H132").Value = 1870.6323
$ws.Range("I132").Value = 1315.9454
$ws.Range("J132").Value = 4217.385
$ws.Range("K132").Value = 3947.8362
$ws.Range("L132").Value = 12652.155
$ws.Range("M132").Value = -1417.8362
$ws.Range("N132").Value = -17712.155

# Row 136 (CRP)
$ws.Range("H136").Value = 2001.4524
$ws.Range("I136").Value = 768.6799999999999
$ws.Range("J136").Value = 3814.353
$ws.Range("K136").Value = 2306.04
$ws.Range("L136").Value = 11443.059
$ws.Range("M136").Value = 243.96
$ws.Range("N136").Value = -16543.059

$ws = $wb.Worksheets.Item("CUL")
# Row 131 (CUL)
$ws.Range("H131").Value = 1638.2894
$ws.Range("I131").Value = 417.66666
$ws.Range("J131").Value = 2017.1034
$ws.Range("K131").Value = 1252.99998
$ws.Range("L131").Value = 6051.3102
$ws.Range("M131").Value = 3787.00002
$ws.Range("N131").Value = -16131.3102

$ws = $wb.Worksheets.Item("GSM")
# Row 57 (GSM)
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

# Row 97 (GSM)
$ws.Range("H97").Value = 948.75
$ws.Range("I97").Value = 787.2222
$ws.Range("J97").Value = 1433.3334
$ws.Range("K97").Value = 787.2222
$ws.Range("L97").Value = 1433.3334
$ws.Range("M97").Value = -291.2222
$ws.Range("N97").Value = -2425.3334

# Row 111 (GSM)
$ws.Range("H111").Value = 31950
$ws.Range("J111").Value = 31950
$ws.Range("L111").Value = 31950
$ws.Range("N111").Value = -38084

# Row 132 (GSM)
$ws.Range("H132").Value = 2295.0334
$ws.Range("I132").Value = 2067.4902
$ws.Range("J132").Value = 3584.4443
$ws.Range("K132").Value = 6202.470600000001
$ws.Range("L132").Value = 10753.3329
$ws.Range("M132").Value = -3672.470600000001
$ws.Range("N132").Value = -15813.3329

# Row 137 (GSM)
$ws.Range("H137").Value = 55000
$ws.Range("J137").Value = 55000
$ws.Range("L137").Value = 55000
$ws.Range("N137").Value = -65200

$ws = $wb.Worksheets.Item("LTW")
# Row 100 (LTW)
$ws.Range("H100").Value = 2571.1428
$ws.Range("I100").Value = 1999
$ws.Range("J100").Value = 2800
$ws.Range("K100").Value = 1999
$ws.Range("L100").Value = 2800
$ws.Range("M100").Value = -1458
$ws.Range("N100").Value = -3882

# Row 110 (LTW)
$ws.Range("H110").Value = 20001
$ws.Range("J110").Value = 20001
$ws.Range("L110").Value = 20001
$ws.Range("N110").Value = -28181

$ws = $wb.Worksheets.Item("WVR")
# Row 58 (WVR)
$ws.Range("H58").Value = 3000
$ws.Range("I58").Value = 3000
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2692
$ws.Range("N58").ClearContents()

# Row 100 (WVR)
$ws.Range("H100").Value = 545.1667
$ws.Range("I100").Value = 331.33334
$ws.Range("J100").Value = 759
$ws.Range("K100").Value = 662.66668
$ws.Range("L100").Value = 1518
$ws.Range("M100").Value = -121.66668
$ws.Range("N100").Value = -2600

# Row 122 (WVR)
$ws.Range("H122").Value = 64280.375
$ws.Range("I122").Value = 92426
$ws.Range("J122").Value = 2360
$ws.Range("K122").Value = 277278
$ws.Range("L122").Value = 7080
$ws.Range("M122").Value = -274828
$ws.Range("N122").Value = -11980
